$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DictionaryPage")

# New row order (Path, ContentType) for rows 2..9 after the edit
$newData = @(
    @("/publications/dictionaries/cancer-terms", "Term Dictionary"),
    @("/espanol/publicaciones/diccionario", "Term Dictionary"),
    @("/publications/dictionaries/cancer-terms?expand=D", "Dictionary Search"),
    @("/publications/dictionaries/cancer-terms/search?contains=true&q=breast", "Dictionary Search"),
    @("/publications/dictionaries/cancer-drug", "Drug Dictionary"),
    @("/publications/dictionaries/genetics-dictionary", "Genetics Dictionary"),
    @("/publications/dictionaries/cancer-terms/def/cadmium", "Definition"),
    @("/espanol/publicaciones/diccionario/def/yin-y-yang", "Definition")
)

$r = 2
foreach ($row in $newData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

# The active sheet moves from CTHPPage to DictionaryPage (tabSelected + activeTab update automatically)
$ws.Activate()
$ws.Range("A10").Select()
